$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 51, shifting existing rows 51-157 down to 52-158.
$ws.Rows(51).Insert()

# Populate the newly inserted row 51 with the new weekly price record.
$ws.Range("A51").Value = 7
$ws.Range("B51").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C51").Value = "Ñuble"
$ws.Range("D51").Value = 44544
$ws.Range("E51").Value = 16
$ws.Range("F51").Value = 100112032
$ws.Range("G51").Value = "Zapallo italiano"
$ws.Range("H51").Value = "Sin especificar"
$ws.Range("I51").Value = "Primera"
$ws.Range("J51").Value = 300
$ws.Range("K51").Value = 5500
$ws.Range("L51").Value = 6000
$ws.Range("M51").Value = 5750
$ws.Range("N51").Value = "`$/caja 60 unidades"
$ws.Range("O51").Value = "Región del Maule"
$ws.Range("P51").Value = 96
$ws.Range("Q51").Value = 60
$ws.Range("R51").Value = "Hortaliza"
